$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.85
$ws.Range("N4").Value = 1.8
$ws.Range("O4").Value = 2
$ws.Range("T4").Value = 8
$ws.Range("X4").Value = 15
$ws.Range("AG4").Value = 13

# Row 9
$ws.Range("G9").Value = 2.25
$ws.Range("H9").Value = 2.7
$ws.Range("I9").Value = 3.8
$ws.Range("Q9").Value = 2.57
$ws.Range("R9").Value = 1.78
$ws.Range("S9").Value = 1.93
$ws.Range("T9").Value = 6.8
$ws.Range("U9").Value = 10.75
$ws.Range("W9").Value = 24
$ws.Range("AA9").Value = 5.3
$ws.Range("AE9").Value = 9.75
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 65

# Row 12
$ws.Range("G12").Value = 3.5
$ws.Range("I12").Value = 1.95
$ws.Range("K12").Value = 10

# Row 17
$ws.Range("G17").Value = 3.4
$ws.Range("T17").Value = 9
$ws.Range("U17").Value = 17.5
$ws.Range("V17").Value = 11.75
$ws.Range("W17").Value = 50
$ws.Range("X17").Value = 32
$ws.Range("Z17").Value = 7.8
$ws.Range("AA17").Value = 6
$ws.Range("AD17").Value = 700
$ws.Range("AE17").Value = 6.4
$ws.Range("AF17").Value = 9.5
$ws.Range("AI17").Value = 19
$ws.Range("AJ17").Value = 35

# Row 18
$ws.Range("H18").Value = 3.75
$ws.Range("M18").Value = 3.25
$ws.Range("N18").Value = 1.72
$ws.Range("O18").Value = 1.88
$ws.Range("T18").Value = 15
$ws.Range("U18").Value = 35
$ws.Range("W18").Value = 110
$ws.Range("Z18").Value = 11
$ws.Range("AA18").Value = 7.4
$ws.Range("AF18").Value = 7.4
$ws.Range("AH18").Value = 11.5

# Row 20
$ws.Range("G20").Value = 2.2
$ws.Range("H20").Value = 3.1
$ws.Range("I20").Value = 3.2
$ws.Range("L20").Value = 1.38
$ws.Range("M20").Value = 2.6
$ws.Range("N20").Value = 2.1
$ws.Range("O20").Value = 1.57
$ws.Range("P20").Value = 1.47
$ws.Range("Q20").Value = 2.32
$ws.Range("R20").Value = 1.85
$ws.Range("S20").Value = 1.75
$ws.Range("T20").Value = 6.6
$ws.Range("U20").Value = 10
$ws.Range("V20").Value = 9
$ws.Range("W20").Value = 21
$ws.Range("X20").Value = 20
$ws.Range("Y20").Value = 35
$ws.Range("Z20").Value = 7.7
$ws.Range("AA20").Value = 6
$ws.Range("AB20").Value = 15.5
$ws.Range("AC20").Value = 90
$ws.Range("AD20").Value = 800
$ws.Range("AE20").Value = 8.5
$ws.Range("AF20").Value = 16
$ws.Range("AG20").Value = 11.5
$ws.Range("AH20").Value = 45
$ws.Range("AI20").Value = 32
$ws.Range("AJ20").Value = 45

# Row 34
$ws.Range("L34").Value = 1.13
$ws.Range("M34").Value = 5.5

# Row 36
$ws.Range("G36").Value = 2.8
$ws.Range("I36").Value = 2.2
$ws.Range("U36").Value = 17
$ws.Range("V36").Value = 11
$ws.Range("W36").Value = 29
$ws.Range("AE36").Value = 9.5
$ws.Range("AF36").Value = 12
$ws.Range("AH36").Value = 21

